$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.046.02"
$ws.Range("E2").Value = "  -2.20%  "
$ws.Range("D3").Value = "1.831.75"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.35"
$ws.Range("E5").Value = "  -2.84%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4659"
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3868"
$ws.Range("E8").Value = "  -1.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07860"
$ws.Range("E9").Value = "  -0.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9593"
$ws.Range("E10").Value = "  -2.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.90"
$ws.Range("E11").Value = "  -1.90%  "
$ws.Range("D12").Value = "1.838.91"
$ws.Range("E12").Value = "  -6.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.686"
$ws.Range("E13").Value = "  -2.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.903"
$ws.Range("E14").Value = "  -1.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06867"
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.22"
$ws.Range("E16").Value = "  -0.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009947"
$ws.Range("E18").Value = "  -1.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.64"
$ws.Range("E19").Value = "  -2.91%  "
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("D21").Value = "28.051.10"
$ws.Range("E21").Value = "  -2.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.324"
$ws.Range("E22").Value = "  -1.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.99"
$ws.Range("E23").Value = "  -2.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.093"
$ws.Range("E24").Value = "  -1.58%  "
$ws.Range("D25").Value = "2.048.24"
$ws.Range("E25").Value = "  -5.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.79"
$ws.Range("E26").Value = "  +0.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.14"
$ws.Range("E27").Value = "  -1.74%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.688"
$ws.Range("E28").Value = "  -6.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.960"
$ws.Range("E29").Value = "  -3.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.65"
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09256"
$ws.Range("E31").Value = "  -1.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9347"
$ws.Range("E32").Value = "  -4.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.270"
$ws.Range("E33").Value = "  -1.88%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.321"
$ws.Range("E34").Value = "  -2.53%  "
$ws.Range("E35").Value = "  -5.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05832"
$ws.Range("E36").Value = "  -5.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02124"
$ws.Range("E37").Value = "  -3.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.140"
$ws.Range("E38").Value = "  -2.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.796"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5597"
$ws.Range("E40").Value = "  -2.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.899"
$ws.Range("E41").Value = "  -2.63%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1762"
$ws.Range("E42").Value = "  -2.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.07192"
$ws.Range("E43").Value = "  +0.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.61"
$ws.Range("E44").Value = "  -2.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5267"
$ws.Range("E45").Value = "  -2.62%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.111"
$ws.Range("E46").Value = "  -11.89%  "
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.115"
$ws.Range("E47").Value = "  -10.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.832"
$ws.Range("E48").Value = "  -4.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "112.39"
$ws.Range("E49").Value = "  -1.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.000"
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("E51").Value = "  +0.19%  "
